$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header cell F1 - copy formatting from E1 (same header style), then set text
$ws.Cells.Item(1, 5).Copy() | Out-Null
$ws.Cells.Item(1, 6).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1, 6).Value = "time_taken"

# Data cells F2:F20 - plain text timestamps
$ws.Cells.Item(2, 6).Value = "2021-10-05 10:52:31.476127"
$ws.Cells.Item(3, 6).Value = "2021-10-05 10:52:31.476137"
$ws.Cells.Item(4, 6).Value = "2021-10-05 10:52:31.476141"
$ws.Cells.Item(5, 6).Value = "2021-10-05 10:52:31.476143"
$ws.Cells.Item(6, 6).Value = "2021-10-05 10:52:31.476146"
$ws.Cells.Item(7, 6).Value = "2021-10-05 10:52:31.476149"
$ws.Cells.Item(8, 6).Value = "2021-10-05 10:52:31.476151"
$ws.Cells.Item(9, 6).Value = "2021-10-05 10:52:31.476154"
$ws.Cells.Item(10, 6).Value = "2021-10-05 10:52:31.476157"
$ws.Cells.Item(11, 6).Value = "2021-10-05 10:52:31.476159"
$ws.Cells.Item(12, 6).Value = "2021-10-05 10:52:31.476162"
$ws.Cells.Item(13, 6).Value = "2021-10-05 10:52:31.476164"
$ws.Cells.Item(14, 6).Value = "2021-10-05 10:52:31.476167"
$ws.Cells.Item(15, 6).Value = "2021-10-05 10:52:31.476170"
$ws.Cells.Item(16, 6).Value = "2021-10-05 10:52:31.476172"
$ws.Cells.Item(17, 6).Value = "2021-10-05 10:52:31.476175"
$ws.Cells.Item(18, 6).Value = "2021-10-05 10:52:31.476178"
$ws.Cells.Item(19, 6).Value = "2021-10-05 10:52:31.476181"
$ws.Cells.Item(20, 6).Value = "2021-10-05 10:52:31.476183"
